# "minor update to delay figure"
# Slide 9 ("Rounded Rectangle 94" shape) currently reads:
#   "The " + "arrival time between packets on node 2 is the inter-arrival time (IAT)"
# It needs to become two runs split differently, with the tail bolded:
#   "The arrival time between packets on node 2 is the " (unchanged formatting)
#   "inter-arrival time (IAT)" (bold)

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(9)
$shp = $s.Shapes.Item(2)
$tr  = $shp.TextFrame.TextRange

$phrase = "inter-arrival time (IAT)"

$full = $tr.Text
$splitAt0 = $full.IndexOf($phrase)   # 0-based index where the bold phrase begins

if ($splitAt0 -ge 0) {
    # Re-flow the leading text into a single run (this merges the original
    # "The " / "arrival time ... " runs) while leaving the trailing phrase,
    # and the paragraph's endParaRPr, untouched.
    $prefix = $full.Substring(0, $splitAt0)
    $leading = $tr.Characters(1, $splitAt0)
    $leading.Delete()
    $null = $tr.InsertBefore($prefix)

    # Recompute the phrase's position now that the prefix run was rebuilt,
    # then bold just that trailing run.
    $current = $tr.Text
    $boldAt0 = $current.IndexOf($phrase)
    $boldRange = $tr.Characters($boldAt0 + 1, $phrase.Length)
    $boldRange.Font.Bold = 1

    # Also try to carry the bold weight onto the paragraph end mark so that
    # any subsequently typed text keeps matching formatting.
    $endMark = $tr.Characters($tr.Length + 1, 1)
    $endMark.Font.Bold = 1
}
